$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 114.056483
$ws.Range("H2").Value = 342.169449
$ws.Range("I2").Value = 0.05607005241237744
$ws.Range("J2").Value = 0.05607005241237744
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2614153333333333
$ws.Range("N2").Value = 0.784246
$ws.Range("Q2").Value = 29.81611352227267
$ws.Range("R2").Value = 268.345021700454
$ws.Range("S2").Value = 0.05607005241237744
$ws.Range("T2").Value = 0.05607005241237744

# Row 3
$ws.Range("I3").Value = 0.1237492078808041
$ws.Range("J3").Value = 0.1237492078808041
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.2614153333333333
$ws.Range("N3").Value = 0.784246
$ws.Range("Q3").Value = 65.80554630712044
$ws.Range("R3").Value = 592.249916764084
$ws.Range("S3").Value = 0.1237492078808041
$ws.Range("T3").Value = 0.1237492078808041

# Row 4
$ws.Range("G4").Value = 70.67310566666667
$ws.Range("H4").Value = 212.019317
$ws.Range("I4").Value = 0.03474282771699605
$ws.Range("J4").Value = 0.03474282771699606
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2614153333333333
$ws.Range("N4").Value = 0.784246
$ws.Range("Q4").Value = 18.47503347555356
$ws.Range("R4").Value = 166.275301279982
$ws.Range("S4").Value = 0.03474282771699605
$ws.Range("T4").Value = 0.03474282771699606

# Row 5
$ws.Range("G5").Value = 1597.720744
$ws.Range("H5").Value = 4793.162232
$ws.Range("I5").Value = 0.7854379119898224
$ws.Range("J5").Value = 0.7854379119898224
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2614153333333333
$ws.Range("N5").Value = 0.784246
$ws.Range("Q5").Value = 417.6687008663413
$ws.Range("R5").Value = 3759.018307797072
$ws.Range("S5").Value = 0.7854379119898224
$ws.Range("T5").Value = 0.7854379119898224
